$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Date, C,D,E,F,G,H,I,J (column B is intentionally left blank,
# matching the existing data rows in the sheet).
$newRows = @(
    @("2024-09-02", 1964.5, 1806.650024414062, 1050.949951171875, 1766.300048828125, 1749.5, 38025.74987792969, 0, 174.4011395709637),
    @("2024-09-03", 1941.25, 1790.449951171875, 1068.800048828125, 1769.300048828125, 1718.75, 37895.55029296875, -0.003423984678248408, 173.8039927412036),
    @("2024-09-04", 1922.449951171875, 1785.25, 1056.199951171875, 1749.699951171875, 1729.550048828125, 37618.74926757812, -0.007304314708473397, 172.5344736806327),
    @("2024-09-05", 1933.150024414062, 1790.550048828125, 1074.900024414062, 1722.900024414062, 1720.75, 37761.05065917969, 0.003782725220059497, 173.1871241855541),
    @("2024-09-06", 1901.849975585938, 1756.099975585938, 1112.650024414062, 1730.300048828125, 1715, 37713.70007324219, -0.001253953084220899, 172.9699556570342),
    @("2024-09-09", 1894.650024414062, 1746.75, 1077.550048828125, 1750.400024414062, 1741.199951171875, 37455.65051269531, -0.006842329446480399, 171.7864382360857),
    @("2024-09-10", 1912.300048828125, 1779.099975585938, 1091, 1756.349975585938, 1745.150024414062, 37856.65014648438, 0.01070598503296977, 173.6255812727084),
    @("2024-09-11", 1910.150024414062, 1778.75, 1077.849975585938, 1789.349975585938, 1782.650024414062, 37910.79992675781, 0.001430390170918655, 173.873933597581),
    @("2024-09-12", 1950.449951171875, 1807.599975585938, 1083.75, 1838.050048828125, 1812.75, 38550.34973144531, 0.01686985782212682, 176.8071621363461),
    @("2024-09-13", 1944.099975585938, 1812.800048828125, 1089.699951171875, 1826.050048828125, 1814.099975585938, 38552.5498046875, 0.00005707012407186833, 176.817252543026),
    @("2024-09-16", 1950.25, 1811.849975585938, 1094.650024414062, 1757.849975585938, 1797.199951171875, 38385.89990234375, -0.004322668751821118, 176.0529301306754),
    @("2024-09-17", 1952.550048828125, 1813.75, 1080.300048828125, 1741.150024414062, 1848.5, 38346.90075683594, -0.001015975803798501, 175.8740646134748),
    @("2024-09-18", 1892.150024414062, 1756.5, 1065.800048828125, 1727.25, 1805.599975585938, 37454.75048828125, -0.023265250931541, 171.7823103678923),
    @("2024-09-19", 1894.199951171875, 1736.5, 1060.75, 1676.449951171875, 1877.449951171875, 37317.94946289062, -0.003652434567236724, 171.1548867194648),
    @("2024-09-20", 1905.75, 1760.050048828125, 1114.699951171875, 1662, 1931.449951171875, 38001.24975585938, 0.01831023147850691, 174.2887723139758),
    @("2024-09-23", 1896.449951171875, 1752.800048828125, 1106.699951171875, 1692.900024414062, 1920.400024414062, 37915.79968261719, -0.002248612184892999, 173.8968644568606),
    @("2024-09-24", 1898.599975585938, 1775.599975585938, 1098.5, 1660.900024414062, 1838.75, 37717.7998046875, -0.00522209420840627, 172.9887586481204),
    @("2024-09-25", 1895.300048828125, 1782.400024414062, 1088.599975585938, 1654.75, 1722.050048828125, 37400.95031738281, -0.008400529430280024, 171.5355614899892),
    @("2024-09-26", 1900.25, 1783.849975585938, 1094.949951171875, 1682.449951171875, 1685.900024414062, 37499.49938964844, 0.002634934979708856, 171.9875465412232),
    @("2024-09-27", 1906.75, 1808.400024414062, 1121, 1671, 1640.800048828125, 37745.10021972656, 0.006549442901254356, 173.1139691570218)
)

$startRow = 657
$r = $startRow
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = "'" + $row[0]   # A: Date, forced text (no auto date conversion)
    $ws.Cells.Item($r, 3).Value = $row[1]          # C
    $ws.Cells.Item($r, 4).Value = $row[2]          # D
    $ws.Cells.Item($r, 5).Value = $row[3]          # E
    $ws.Cells.Item($r, 6).Value = $row[4]          # F
    $ws.Cells.Item($r, 7).Value = $row[5]          # G
    $ws.Cells.Item($r, 8).Value = $row[6]          # H
    $ws.Cells.Item($r, 9).Value = $row[7]          # I
    $ws.Cells.Item($r, 10).Value = $row[8]         # J
    $r++
}

$endRow = $r - 1
Write-Output "Added rows $startRow to $endRow"
